# Updates the StructureDefinition-episode-body-system workbook so that it reflects
# the new LinuxForHealth publication (regenerated from the FHIR IG Publisher at
# version 8.0.0) instead of the old Alvearie / IBM 7.0.0 publication.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-body-system"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The Fixed Value for Extension.url carries the same canonical URL as the
# Metadata URL and must be refreshed to match.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-body-system"

# The root "Extension" row no longer repeats the ele-1/ext-1 constraint text
# in the regenerated IG output; it is now cleared.
$elements.Range("AI2").Value = ""
